# The deck's slide-master theme (theme1.xml) is being re-coloured from the
# "Integral" palette to the "Office Theme" palette (the palette that used to
# live in theme2.xml, the notes-master theme). The font scheme and format
# scheme are already identical between the two themes, so only the 10
# colour-scheme entries that actually differ need to change - dk1/lt1 are
# black/white in both so they're untouched.
#
# PowerPoint's theme colour scheme exposes exactly this: 12 theme colours in
# clrScheme document order -> dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink,
# addressed as ThemeColorScheme.Colors(index).RGB (an OLE "BGR" long). Going
# through Slide.ThemeColorScheme (rather than SlideMaster.ColorScheme)
# updates the same underlying theme1.xml colours for every slide (there is
# only the one slide master/theme for the deck) while leaving the rest of
# the theme part - names, fonts, fill/line/effect styles - untouched.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

$cs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      -> 44546A
$cs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  -> ED7D31
$cs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  -> FFC000
$cs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  -> 4472C4
$cs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  -> 70AD47
$cs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    -> 0563C1
$cs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink -> 954F72
